$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 688.8946999999999
$ws.Range("I15").Value = 688.8946999999999
$ws.Range("K15").Value = 2066.6841
$ws.Range("M15").Value = -1897.6841
$ws.Range("H17").Value = 1500
$ws.Range("J17").Value = 1500
$ws.Range("L17").Value = 4500
$ws.Range("N17").Value = -4836
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 3500
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 3500
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -4436
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 3500
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 3500
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -3968
$ws.Range("H29").Value = 4549.8
$ws.Range("J29").Value = 5637.25
$ws.Range("L29").Value = 16911.75
$ws.Range("N29").Value = -17473.75
$ws.Range("H70").Value = 3832.9512
$ws.Range("I70").Value = 1301.7333
$ws.Range("K70").Value = 3905.199900000001
$ws.Range("M70").Value = -3635.199900000001
$ws.Range("H73").Value = 3832.9512
$ws.Range("I73").Value = 1301.7333
$ws.Range("K73").Value = 3905.199900000001
$ws.Range("M73").Value = -2969.199900000001
$ws.Range("H106").Value = 8930.4
$ws.Range("I106").Value = 8533.888999999999
$ws.Range("K106").Value = 8533.888999999999
$ws.Range("M106").Value = -7902.888999999999
$ws.Range("H107").Value = 1003.06665
$ws.Range("I107").Value = 1548.2222
$ws.Range("K107").Value = 1548.2222
$ws.Range("M107").Value = 371.7778000000001
$ws.Range("H115").Value = 433.5
$ws.Range("I115").Value = 236.77777
$ws.Range("K115").Value = 710.33331
$ws.Range("M115").Value = 856.66669
$ws.Range("H116").Value = 7178.2
$ws.Range("I116").Value = 15599.333
$ws.Range("J116").Value = 3569.1428
$ws.Range("K116").Value = 15599.333
$ws.Range("L116").Value = 3569.1428
$ws.Range("M116").Value = -12157.333
$ws.Range("N116").Value = -10453.1428
$ws.Range("H129").Value = 1005.75
$ws.Range("I129").Value = 918.58826
$ws.Range("J129").Value = 1499.6666
$ws.Range("K129").Value = 2755.76478
$ws.Range("L129").Value = 4498.9998
$ws.Range("M129").Value = 2244.23522
$ws.Range("N129").Value = -14498.9998
$ws.Range("H132").Value = 3005.7334
$ws.Range("I132").Value = 2851.1428
$ws.Range("K132").Value = 8553.428400000001
$ws.Range("M132").Value = -6023.428400000001
$ws.Range("H135").Value = 3185.182
$ws.Range("J135").Value = 2916
$ws.Range("L135").Value = 26244
$ws.Range("N135").Value = -31314
$ws.Range("H138").Value = 4287.325
$ws.Range("I138").Value = 4618.25
$ws.Range("J138").Value = 4066.7083
$ws.Range("K138").Value = 13854.75
$ws.Range("L138").Value = 12200.1249
$ws.Range("M138").Value = -8714.75
$ws.Range("N138").Value = -22480.1249
$ws.Range("H141").Value = 7391.625
$ws.Range("J141").Value = 29999
$ws.Range("L141").Value = 89997
$ws.Range("N141").Value = -100357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1593
$ws.Range("I110").Value = 1593
$ws.Range("K110").Value = 1593
$ws.Range("M110").Value = 452
$ws.Range("H139").Value = 56885.75
$ws.Range("J139").Value = 56885.75
$ws.Range("L139").Value = 56885.75
$ws.Range("N139").Value = -67165.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1598.96
$ws.Range("I20").Value = 1738.5555
$ws.Range("K20").Value = 1738.5555
$ws.Range("M20").Value = -1491.5555
$ws.Range("H99").Value = 4161.273
$ws.Range("J99").Value = 4862.2
$ws.Range("L99").Value = 4862.2
$ws.Range("N99").Value = -7858.2
$ws.Range("H107").Value = 18524.375
$ws.Range("I107").Value = 20442.143
$ws.Range("K107").Value = 20442.143
$ws.Range("M107").Value = -18522.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2119.1538
$ws.Range("I58").Value = 1233.7142
$ws.Range("J58").Value = 3152.1667
$ws.Range("K58").Value = 1233.7142
$ws.Range("L58").Value = 3152.1667
$ws.Range("M58").Value = -1030.7142
$ws.Range("N58").Value = -3558.1667
$ws.Range("H136").Value = 2119.1538
$ws.Range("I136").Value = 1233.7142
$ws.Range("J136").Value = 3152.1667
$ws.Range("K136").Value = 3701.1426
$ws.Range("L136").Value = 9456.500100000001
$ws.Range("M136").Value = -1151.1426
$ws.Range("N136").Value = -14556.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 370.05264
$ws.Range("I5").Value = 347.6875
$ws.Range("J5").Value = 489.33334
$ws.Range("K5").Value = 1043.0625
$ws.Range("L5").Value = 1468.00002
$ws.Range("M5").Value = -931.0625
$ws.Range("N5").Value = -1692.00002
$ws.Range("H95").Value = 8500
$ws.Range("J95").Value = 8500
$ws.Range("L95").Value = 25500
$ws.Range("N95").Value = -29618
$ws.Range("H107").Value = 642.46155
$ws.Range("J107").Value = 1356.2
$ws.Range("L107").Value = 4068.6
$ws.Range("N107").Value = -7908.6
$ws.Range("H135").Value = 370.05264
$ws.Range("I135").Value = 347.6875
$ws.Range("J135").Value = 489.33334
$ws.Range("K135").Value = 3129.1875
$ws.Range("L135").Value = 4404.00006
$ws.Range("M135").Value = -594.1875
$ws.Range("N135").Value = -9474.00006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 217.27272
$ws.Range("J2").Value = 599
$ws.Range("L2").Value = 599
$ws.Range("N2").Value = -825
$ws.Range("H3").Value = 4862.5
$ws.Range("I3").Value = 1850.25
$ws.Range("J3").Value = 7874.75
$ws.Range("K3").Value = 1850.25
$ws.Range("L3").Value = 7874.75
$ws.Range("M3").Value = -1734.25
$ws.Range("N3").Value = -8106.75
$ws.Range("H97").Value = 1639.6957
$ws.Range("I97").Value = 523.53845
$ws.Range("K97").Value = 523.53845
$ws.Range("M97").Value = -27.53845000000001
$ws.Range("H102").Value = 12401.296
$ws.Range("I102").Value = 22952.75
$ws.Range("K102").Value = 22952.75
$ws.Range("M102").Value = -21330.75
$ws.Range("H122").Value = 3103.6667
$ws.Range("I122").Value = 2867.647
$ws.Range("K122").Value = 8602.940999999999
$ws.Range("M122").Value = -6152.940999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1000.7
$ws.Range("I55").Value = 134.5
$ws.Range("J55").Value = 2300
$ws.Range("K55").Value = 134.5
$ws.Range("L55").Value = 2300
$ws.Range("M55").Value = 38.5
$ws.Range("N55").Value = -2646
$ws.Range("H132").Value = 2181.5557
$ws.Range("I132").Value = 2010.1666
$ws.Range("K132").Value = 6030.4998
$ws.Range("M132").Value = -3500.4998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 419.14285
$ws.Range("I100").Value = 429.07693
$ws.Range("K100").Value = 858.15386
$ws.Range("M100").Value = -317.15386
$ws.Range("H132").Value = 4907.1377
$ws.Range("I132").Value = 4721.5415
$ws.Range("K132").Value = 14164.6245
$ws.Range("M132").Value = -11634.6245
